$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WIP")

# Remove the "t-shirt size" column (column F) entirely.
$ws.Columns.Item(6).Delete()

# Add a new row of data (row 23) for the new task.
$ws.Range("A23").Value = "RMIT final report"
$ws.Range("B23").Value = "Don't even know what I'm doing"
$ws.Range("C23").Value = "Thao"
$ws.Range("D23").Value = "Done"

# Match the style of the row above (A22:D22 -> Danh project row) for the new row.
$ws.Range("A22:D22").Copy()
$ws.Range("A23:D23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to mirror the authored workbook state.
$ws.Range("E26").Select() | Out-Null
